# Create cinematics at boss spawn
# Update SpawnData sheet: fix numbers for existing rows, convert the old
# "template/example" rows (4-6) into real spawn data rows, add new rows
# for additional spawns (7-11), and append a run of index-only rows (12-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$monsterPath = '("/Game/Character/Monster/WhiteMinion/BWhiteMinion_BP.BWhiteMinion_BP_C")'

# --- Fix existing data rows 2 and 3 ---
$ws.Range("D2").Value = -1
$ws.Range("F2").Value = -300
$ws.Range("H2").Value = -1

$ws.Range("C3").Value = -27
$ws.Range("D3").Value = -1
$ws.Range("E3").Value = -10
$ws.Range("F3").Value = -300
$ws.Range("H3").Value = -1

# --- Replace rows 4-6 (previously placeholder/template text rows) and
#     add new rows 7-11 with real spawn entries ---
$rowsData = @(
    @{ Row = 4;  A = 3;  C = -90 },
    @{ Row = 5;  A = 4;  C = -90 },
    @{ Row = 6;  A = 5;  C = -9  },
    @{ Row = 7;  A = 6;  C = -9  },
    @{ Row = 8;  A = 7;  C = -9  },
    @{ Row = 9;  A = 8;  C = -9  },
    @{ Row = 10; A = 9;  C = -9  },
    @{ Row = 11; A = 10; C = -9  }
)

foreach ($entry in $rowsData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $monsterPath
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = -1
    $ws.Cells.Item($r, 5).Value = -10
    $ws.Cells.Item($r, 6).Value = -300
    $ws.Cells.Item($r, 7).Value = -1
    $ws.Cells.Item($r, 8).Value = -1
}

# --- Append trailing index-only rows 12-24 (column A = row-1) ---
for ($r = 12; $r -le 24; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# --- Re-apply formatting (number format / alignment) that plain .Value
#     assignment does not carry over, by copying the existing formats
#     from rows 2/3 onto the new / rewritten cells. ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4:A24").PasteSpecial(-4122) | Out-Null

$ws.Range("C2:H2").Copy() | Out-Null
$ws.Range("C4:H11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Column widths / view tweaks ---
# (engine quantizes stored width to whole "pixels"; 88 is the closest
# achievable ColumnWidth to the authored 88.75 character width)
$ws.Columns.Item(2).ColumnWidth = 88

$excel.ActiveWindow.ScrollRow = 5
$ws.Range("B18").Select() | Out-Null
